# Update "想去人数" (want-to-go count) values in column F for the
# "展览" sheet and the consolidated "全部类型" sheet, matching the
# refreshed data snapshot published at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 168
$ws1.Range("F3").Value  = 1405
$ws1.Range("F7").Value  = 524
$ws1.Range("F8").Value  = 856
$ws1.Range("F9").Value  = 576
$ws1.Range("F12").Value = 526
$ws1.Range("F13").Value = 101
$ws1.Range("F15").Value = 539
$ws1.Range("F17").Value = 442
$ws1.Range("F19").Value = 273
$ws1.Range("F20").Value = 40
$ws1.Range("F22").Value = 512
$ws1.Range("F23").Value = 497
$ws1.Range("F24").Value = 5
$ws1.Range("F25").Value = 413

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 168
$ws4.Range("F4").Value  = 1405
$ws4.Range("F12").Value = 524
$ws4.Range("F13").Value = 856
$ws4.Range("F14").Value = 576
$ws4.Range("F17").Value = 526
$ws4.Range("F18").Value = 101
$ws4.Range("F20").Value = 539
$ws4.Range("F24").Value = 442
$ws4.Range("F28").Value = 273
$ws4.Range("F29").Value = 40
$ws4.Range("F33").Value = 512
$ws4.Range("F36").Value = 497
$ws4.Range("F37").Value = 5
$ws4.Range("F38").Value = 413

$wb.Save()
